# Add review comments for a few still-unresolved bugs.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Bug in row 20 (register page, no avatar after success) - explain slow network cause.
$ws.Cells.Item(20, 3).Value = "这个是由于网络慢的问题。注册时会发两个api调用，其中第一个是传文本注册信息，第二个是上传photo。第二个会比较慢，失败的概率在网慢的时候较大，从而导致这种情况。"

# Bug in row 30 ("is saved" popup on message input) - append follow-up note.
$ws.Cells.Item(30, 3).Value = "拿手机中的message进行测试，test文本在键盘给出的选项中选中后也会出现“输入文字+is saved”的字样，联宇可以作为参考。`n这么看来，不是我们app的问题了，而是那个htc手机系统本身的问题。"

# Bug in row 23 (slow response on edit/preview buttons) - can't be fixed currently.
$ws.Cells.Item(23, 3).Value = "这个目前改不了。"
